$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.161.63"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").Value = "1.656.45"
$ws.Range("E3").Value = "  -0.10%  "
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5241"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.32%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("E8").Value = "  +1.50%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06356"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.57"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07692"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.49%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.626"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.95%  "
$ws.Range("D13").Value = "1.678.39"
$ws.Range("E13").Value = "  +1.19%  "
$ws.Range("D14").Value = "1.884.34"
$ws.Range("E14").Value = "  -0.10%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5620"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.29%  "
$ws.Range("D16").Value = "0.0₅8203"
$ws.Range("E16").Value = "  +2.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.56"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.70%  "
$ws.Range("D18").Value = "26.149.48"
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("E19").Value = "  -0.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.662"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "193.18"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.957"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("E24").Value = "  -0.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.26"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1198"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.264"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.96"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.24%  "
$ws.Range("E29").Value = "  +2.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05490"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.275"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.468"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.365"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.37%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.564"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9523"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.34%  "
$ws.Range("E36").Value = "  -0.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.402"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5689"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01589"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.875"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.30%  "
$ws.Range("E41").Value = "  -0.21%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "1.028.00"
$ws.Range("E42").Value = "  -3.52%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8303"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.16"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.56%  "
$ws.Range("D45").Value = "1.795.34"
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "58.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₈105"
$ws.Range("E47").Value = "  +2.39%  "
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9999"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.050"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.61%  "
$ws.Range("E50").Value = "  -1.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05210"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.21%  "
